# Update "想去人数" (F column) counts across sheets, reflecting the
# regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 8111
$ws1.Range("F4").Value  = 1910
$ws1.Range("F5").Value  = 6497
$ws1.Range("F6").Value  = 158
$ws1.Range("F7").Value  = 2048
$ws1.Range("F8").Value  = 564
$ws1.Range("F10").Value = 19
$ws1.Range("F15").Value = 8465
$ws1.Range("F16").Value = 158
$ws1.Range("F17").Value = 63
$ws1.Range("F19").Value = 114
$ws1.Range("F20").Value = 1804
$ws1.Range("F21").Value = 858
$ws1.Range("F22").Value = 12
$ws1.Range("F25").Value = 19
$ws1.Range("F27").Value = 184
$ws1.Range("F28").Value = 4
$ws1.Range("F30").Value = 2051
$ws1.Range("F32").Value = 463
$ws1.Range("F34").Value = 9
$ws1.Range("F35").Value = 168
$ws1.Range("F37").Value = 13
$ws1.Range("F39").Value = 3961

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value  = 387
$ws2.Range("F3").Value  = 207
$ws2.Range("F16").Value = 95
$ws2.Range("F20").Value = 6

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 2318

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 2318
$ws4.Range("F4").Value  = 387
$ws4.Range("F6").Value  = 8111
$ws4.Range("F9").Value  = 1910
$ws4.Range("F10").Value = 6497
$ws4.Range("F11").Value = 2048
$ws4.Range("F13").Value = 564
$ws4.Range("F15").Value = 19
$ws4.Range("F23").Value = 8465
$ws4.Range("F24").Value = 158
$ws4.Range("F25").Value = 63
$ws4.Range("F27").Value = 114
$ws4.Range("F28").Value = 1804
$ws4.Range("F29").Value = 858
$ws4.Range("F30").Value = 12
$ws4.Range("F33").Value = 184
$ws4.Range("F35").Value = 2051
$ws4.Range("F38").Value = 463
$ws4.Range("F43").Value = 6
$ws4.Range("F44").Value = 3961

$wb.Save()
